$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10: "blog ser: 179" -> "blog ser: 175"
$ws.Range("B10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 175"

# F10: meetup date changes from 2020,6,5,... to 2023,2,2,...
$ws.Range("F10").Value = "type: meetup`nwidth: 2`nheight: 1`nh3: Meetup coming in`ndate: 2023,2,2,10,30,0,0`nbutton.default: Speak*goto(`"https://forms.gle/dyydXFRSsKzeH4hZ6`")`nbutton.default: Attend*goto(`"https://youtu.be/vscn-HP932E`")`nbutton.default: Details*goto(`"https://www.meetup.com/techshek/events/270179438/`")"

# Update the active cell selection to F10 (matches author's final selection state)
$ws.Range("F10").Select()
